$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Remove the now-obsolete trailing rows (48-51) from Recommandations sheet
$ws1.Range("A48:G51").EntireRow.Delete()

# Update Recommandations sheet data rows (2-47)
$ws1.Range("A2").Value = "BRVM - SERVICES PUBLICS"
$ws1.Range("B2").Value = 0
$ws1.Range("C2").Value = 8
$ws1.Range("D2").Value = 3433.99
$ws1.Range("E2").Value = 113.05
$ws1.Range("F2").Value = "🟡 Observer"
$ws1.Range("G2").Value = "➖ Neutre"

$ws1.Range("A3").Value = "SUCRIVOIRE"
$ws1.Range("B3").Value = 0
$ws1.Range("C3").Value = 3
$ws1.Range("D3").Value = 2945
$ws1.Range("E3").Value = 985
$ws1.Range("F3").Value = "🟡 Observer"
$ws1.Range("G3").Value = "➖ Neutre"

$ws1.Range("A4").Value = "SAFCA CI"
$ws1.Range("B4").Value = 0
$ws1.Range("C4").Value = 4
$ws1.Range("D4").Value = 2760
$ws1.Range("E4").Value = 690
$ws1.Range("F4").Value = "🟡 Observer"
$ws1.Range("G4").Value = "➖ Neutre"

$ws1.Range("A5").Value = "CFAO MOTORS CI"
$ws1.Range("B5").Value = 0
$ws1.Range("C5").Value = 4
$ws1.Range("D5").Value = 2710
$ws1.Range("E5").Value = 675
$ws1.Range("F5").Value = "🟡 Observer"
$ws1.Range("G5").Value = "➖ Neutre"

$ws1.Range("A6").Value = "BRVM - AUTRES SECTEURS"
$ws1.Range("B6").Value = 0
$ws1.Range("C6").Value = 4
$ws1.Range("D6").Value = 2655.54
$ws1.Range("E6").Value = 665.88
$ws1.Range("F6").Value = "🟡 Observer"
$ws1.Range("G6").Value = "➖ Neutre"

$ws1.Range("A7").Value = "NEI-CEDA CI"
$ws1.Range("B7").Value = 0
$ws1.Range("C7").Value = 4
$ws1.Range("D7").Value = 2370
$ws1.Range("E7").Value = 595
$ws1.Range("F7").Value = "🟡 Observer"
$ws1.Range("G7").Value = "➖ Neutre"

$ws1.Range("A8").Value = "UNIWAX CI"
$ws1.Range("B8").Value = 0
$ws1.Range("C8").Value = 4
$ws1.Range("D8").Value = 2370
$ws1.Range("E8").Value = 600
$ws1.Range("F8").Value = "🟡 Observer"
$ws1.Range("G8").Value = "➖ Neutre"

$ws1.Range("A9").Value = "SETAO CI"
$ws1.Range("B9").Value = 0
$ws1.Range("C9").Value = 4
$ws1.Range("D9").Value = 2180
$ws1.Range("E9").Value = 575
$ws1.Range("F9").Value = "🟡 Observer"
$ws1.Range("G9").Value = "➖ Neutre"

$ws1.Range("A10").Value = "AIR LIQUIDE CI"
$ws1.Range("B10").Value = 0
$ws1.Range("C10").Value = 4
$ws1.Range("D10").Value = 2165
$ws1.Range("E10").Value = 540
$ws1.Range("F10").Value = "🟡 Observer"
$ws1.Range("G10").Value = "➖ Neutre"

$ws1.Range("A11").Value = "BRVM - DISTRIBUTION"
$ws1.Range("B11").Value = 0
$ws1.Range("C11").Value = 4
$ws1.Range("D11").Value = 1471.9
$ws1.Range("E11").Value = 370.28
$ws1.Range("F11").Value = "🟡 Observer"
$ws1.Range("G11").Value = "➖ Neutre"

$ws1.Range("A12").Value = "BRVM - TRANSPORT"
$ws1.Range("B12").Value = 0
$ws1.Range("C12").Value = 4
$ws1.Range("D12").Value = 1461.01
$ws1.Range("E12").Value = 348.8
$ws1.Range("F12").Value = "🟡 Observer"
$ws1.Range("G12").Value = "➖ Neutre"

$ws1.Range("A13").Value = "BRVM - AGRICULTURE"
$ws1.Range("B13").Value = 0
$ws1.Range("C13").Value = 4
$ws1.Range("D13").Value = 1327.88
$ws1.Range("E13").Value = 333.59
$ws1.Range("F13").Value = "🟡 Observer"
$ws1.Range("G13").Value = "➖ Neutre"

$ws1.Range("A14").Value = "BRVM - INDUSTRIE"
$ws1.Range("B14").Value = 0
$ws1.Range("C14").Value = 4
$ws1.Range("D14").Value = 773.06
$ws1.Range("E14").Value = 192.93
$ws1.Range("F14").Value = "🟡 Observer"
$ws1.Range("G14").Value = "➖ Neutre"

$ws1.Range("A15").Value = "BRVM-PRINCIPAL"
$ws1.Range("B15").Value = 0
$ws1.Range("C15").Value = 4
$ws1.Range("D15").Value = 709.66
$ws1.Range("E15").Value = 176.76
$ws1.Range("F15").Value = "🟡 Observer"
$ws1.Range("G15").Value = "➖ Neutre"

$ws1.Range("A16").Value = "BRVM - CONSOMMATION DE BASE"
$ws1.Range("B16").Value = 0
$ws1.Range("C16").Value = 4
$ws1.Range("D16").Value = 682.83
$ws1.Range("E16").Value = 171.42
$ws1.Range("F16").Value = "🟡 Observer"
$ws1.Range("G16").Value = "➖ Neutre"

$ws1.Range("A17").Value = "BRVM - INDUSTRIELS"
$ws1.Range("B17").Value = 0
$ws1.Range("C17").Value = 4
$ws1.Range("D17").Value = 549.01
$ws1.Range("E17").Value = 129.78
$ws1.Range("F17").Value = "🟡 Observer"
$ws1.Range("G17").Value = "➖ Neutre"

$ws1.Range("A18").Value = "BRVM-PRESTIGE"
$ws1.Range("B18").Value = 0
$ws1.Range("C18").Value = 4
$ws1.Range("D18").Value = 525.8
$ws1.Range("E18").Value = 131.19
$ws1.Range("F18").Value = "🟡 Observer"
$ws1.Range("G18").Value = "➖ Neutre"

$ws1.Range("A19").Value = "BRVM - FINANCES"
$ws1.Range("B19").Value = 0
$ws1.Range("C19").Value = 4
$ws1.Range("D19").Value = 493.83
$ws1.Range("E19").Value = 123.18
$ws1.Range("F19").Value = "🟡 Observer"
$ws1.Range("G19").Value = "➖ Neutre"

$ws1.Range("A20").Value = "BRVM - SERVICES FINANCIERS"
$ws1.Range("B20").Value = 0
$ws1.Range("C20").Value = 4
$ws1.Range("D20").Value = 485.33
$ws1.Range("E20").Value = 121.06
$ws1.Range("F20").Value = "🟡 Observer"
$ws1.Range("G20").Value = "➖ Neutre"

$ws1.Range("A21").Value = "BRVM - ENERGIE"
$ws1.Range("B21").Value = 0
$ws1.Range("C21").Value = 4
$ws1.Range("D21").Value = 437.99
$ws1.Range("E21").Value = 110.27
$ws1.Range("F21").Value = "🟡 Observer"
$ws1.Range("G21").Value = "➖ Neutre"

$ws1.Range("A22").Value = "BRVM - CONSOMMATION DISCRETIONNAIRE"
$ws1.Range("B22").Value = 0
$ws1.Range("C22").Value = 4
$ws1.Range("D22").Value = 427.25
$ws1.Range("E22").Value = 107.22
$ws1.Range("F22").Value = "🟡 Observer"
$ws1.Range("G22").Value = "➖ Neutre"

$ws1.Range("A23").Value = "BRVM - TELECOMMUNICATIONS"
$ws1.Range("B23").Value = 0
$ws1.Range("C23").Value = 4
$ws1.Range("D23").Value = 387.67
$ws1.Range("E23").Value = 96.01
$ws1.Range("F23").Value = "🟡 Observer"
$ws1.Range("G23").Value = "➖ Neutre"

$ws1.Range("A24").Value = "UNILEVER CI (UNLC)"
$ws1.Range("B24").Value = 4
$ws1.Range("C24").Value = 0
$ws1.Range("D24").Value = 29.92
$ws1.Range("E24").Value = 7.49
$ws1.Range("F24").Value = "🟢 Achat"
$ws1.Range("G24").Value = "✅ Renforcer"

$ws1.Range("A25").Value = "ONATEL BF (ONTBF)"
$ws1.Range("B25").Value = 1
$ws1.Range("C25").Value = 0
$ws1.Range("D25").Value = 7.48
$ws1.Range("E25").Value = 7.48
$ws1.Range("F25").Value = "🟡 Observer"
$ws1.Range("G25").Value = "➖ Neutre"

$ws1.Range("A26").Value = "SETAO CI (STAC)"
$ws1.Range("B26").Value = 2
$ws1.Range("C26").Value = 2
$ws1.Range("D26").Value = 4.97
$ws1.Range("E26").Value = -2.59
$ws1.Range("F26").Value = "🟡 Observer"
$ws1.Range("G26").Value = "👀 À surveiller"

$ws1.Range("A27").Value = "TRACTAFRIC MOTORS CI (PRSC)"
$ws1.Range("B27").Value = 1
$ws1.Range("C27").Value = 0
$ws1.Range("D27").Value = 4.16
$ws1.Range("E27").Value = 4.16
$ws1.Range("F27").Value = "🟡 Observer"
$ws1.Range("G27").Value = "➖ Neutre"

$ws1.Range("A28").Value = "CIE CI (CIEC)"
$ws1.Range("B28").Value = 1
$ws1.Range("C28").Value = 0
$ws1.Range("D28").Value = 4
$ws1.Range("E28").Value = 4
$ws1.Range("F28").Value = "🟡 Observer"
$ws1.Range("G28").Value = "➖ Neutre"

$ws1.Range("A29").Value = "TOTALENERGIES MARKETING CI (TTLC)"
$ws1.Range("B29").Value = 1
$ws1.Range("C29").Value = 0
$ws1.Range("D29").Value = 3.39
$ws1.Range("E29").Value = 3.39
$ws1.Range("F29").Value = "🟡 Observer"
$ws1.Range("G29").Value = "➖ Neutre"

$ws1.Range("A30").Value = "CFAO MOTORS CI (CFAC)"
$ws1.Range("B30").Value = 1
$ws1.Range("C30").Value = 0
$ws1.Range("D30").Value = 3.03
$ws1.Range("E30").Value = 3.03
$ws1.Range("F30").Value = "🟡 Observer"
$ws1.Range("G30").Value = "➖ Neutre"

$ws1.Range("A31").Value = "NSIA BANQUE COTE D'IVOIRE (NSBC)"
$ws1.Range("B31").Value = 1
$ws1.Range("C31").Value = 0
$ws1.Range("D31").Value = 3.03
$ws1.Range("E31").Value = 3.03
$ws1.Range("F31").Value = "🟡 Observer"
$ws1.Range("G31").Value = "➖ Neutre"

$ws1.Range("A32").Value = "TOTALENERGIES MARKETING SN (TTLS)"
$ws1.Range("B32").Value = 1
$ws1.Range("C32").Value = 0
$ws1.Range("D32").Value = 3.02
$ws1.Range("E32").Value = 3.02
$ws1.Range("F32").Value = "🟡 Observer"
$ws1.Range("G32").Value = "➖ Neutre"

$ws1.Range("A33").Value = "BANK OF AFRICA BN (BOAB)"
$ws1.Range("B33").Value = 1
$ws1.Range("C33").Value = 0
$ws1.Range("D33").Value = 2.81
$ws1.Range("E33").Value = 2.81
$ws1.Range("F33").Value = "🟡 Observer"
$ws1.Range("G33").Value = "➖ Neutre"

$ws1.Range("A34").Value = "TOTAL"
$ws1.Range("B34").Value = 0
$ws1.Range("C34").Value = 4
$ws1.Range("D34").Value = 0
$ws1.Range("E34").Value = 0
$ws1.Range("F34").Value = "🟡 Observer"
$ws1.Range("G34").Value = "➖ Neutre"

$ws1.Range("A35").Value = "ORANGE COTE D'IVOIRE (ORAC)"
$ws1.Range("B35").Value = 1
$ws1.Range("C35").Value = 1
$ws1.Range("D35").Value = -0.56
$ws1.Range("E35").Value = 2.6
$ws1.Range("F35").Value = "🟡 Observer"
$ws1.Range("G35").Value = "👀 À surveiller"

$ws1.Range("A36").Value = "BANK OF AFRICA NG (BOAN)"
$ws1.Range("B36").Value = 1
$ws1.Range("C36").Value = 1
$ws1.Range("D36").Value = -1.24
$ws1.Range("E36").Value = 6
$ws1.Range("F36").Value = "🟡 Observer"
$ws1.Range("G36").Value = "👀 À surveiller"

$ws1.Range("A37").Value = "VIVO ENERGY CI (SHEC)"
$ws1.Range("B37").Value = 1
$ws1.Range("C37").Value = 2
$ws1.Range("D37").Value = -2.57
$ws1.Range("E37").Value = 4.43
$ws1.Range("F37").Value = "🟡 Observer"
$ws1.Range("G37").Value = "👀 À surveiller"

$ws1.Range("A38").Value = "BANK OF AFRICA ML (BOAM)"
$ws1.Range("B38").Value = 0
$ws1.Range("C38").Value = 1
$ws1.Range("D38").Value = -2.93
$ws1.Range("E38").Value = -2.93
$ws1.Range("F38").Value = "🟡 Observer"
$ws1.Range("G38").Value = "➖ Neutre"

$ws1.Range("A39").Value = "SMB CI (SMBC)"
$ws1.Range("B39").Value = 1
$ws1.Range("C39").Value = 1
$ws1.Range("D39").Value = -2.93
$ws1.Range("E39").Value = 3.28
$ws1.Range("F39").Value = "🟡 Observer"
$ws1.Range("G39").Value = "👀 À surveiller"

$ws1.Range("A40").Value = "SODE CI (SDCC)"
$ws1.Range("B40").Value = 1
$ws1.Range("C40").Value = 1
$ws1.Range("D40").Value = -2.97
$ws1.Range("E40").Value = 2.43
$ws1.Range("F40").Value = "🟡 Observer"
$ws1.Range("G40").Value = "👀 À surveiller"

$ws1.Range("A41").Value = "SONATEL SN (SNTS)"
$ws1.Range("B41").Value = 0
$ws1.Range("C41").Value = 1
$ws1.Range("D41").Value = -4.23
$ws1.Range("E41").Value = -4.23
$ws1.Range("F41").Value = "🟡 Observer"
$ws1.Range("G41").Value = "➖ Neutre"

$ws1.Range("A42").Value = "BANK OF AFRICA BF (BOABF)"
$ws1.Range("B42").Value = 0
$ws1.Range("C42").Value = 1
$ws1.Range("D42").Value = -5.33
$ws1.Range("E42").Value = -5.33
$ws1.Range("F42").Value = "🟡 Observer"
$ws1.Range("G42").Value = "➖ Neutre"

$ws1.Range("A43").Value = "ECOBANK TRANS. INCORP. TG (ETIT)"
$ws1.Range("B43").Value = 0
$ws1.Range("C43").Value = 1
$ws1.Range("D43").Value = -5.56
$ws1.Range("E43").Value = -5.56
$ws1.Range("F43").Value = "🟡 Observer"
$ws1.Range("G43").Value = "➖ Neutre"

$ws1.Range("A44").Value = "SOLIBRA CI (SLBC)"
$ws1.Range("B44").Value = 0
$ws1.Range("C44").Value = 1
$ws1.Range("D44").Value = -6.67
$ws1.Range("E44").Value = -6.67
$ws1.Range("F44").Value = "🟡 Observer"
$ws1.Range("G44").Value = "➖ Neutre"

$ws1.Range("A45").Value = "SERVAIR ABIDJAN CI (ABJC)"
$ws1.Range("B45").Value = 1
$ws1.Range("C45").Value = 2
$ws1.Range("D45").Value = -8.38
$ws1.Range("E45").Value = 4.24
$ws1.Range("F45").Value = "🟡 Observer"
$ws1.Range("G45").Value = "👀 À surveiller"

$ws1.Range("A46").Value = "AFRICA GLOBAL LOGISTICS CI (SDSC)"
$ws1.Range("B46").Value = 0
$ws1.Range("C46").Value = 2
$ws1.Range("D46").Value = -11.41
$ws1.Range("E46").Value = -7.42
$ws1.Range("F46").Value = "🟡 Observer"
$ws1.Range("G46").Value = "➖ Neutre"

$ws1.Range("A47").Value = "FILTISAC CI (FTSC)"
$ws1.Range("B47").Value = 0
$ws1.Range("C47").Value = 3
$ws1.Range("D47").Value = -20.52
$ws1.Range("E47").Value = -7.4
$ws1.Range("F47").Value = "🔴 Vente"
$ws1.Range("G47").Value = "⚠️ Risque de décrochage"

# Update Top_YTD sheet data rows (2-11)
$ws2.Range("A2").Value = "BRVM - SERVICES PUBLICS"
$ws2.Range("B2").Value = 10644402.42
$ws2.Range("A3").Value = "SAFCA CI"
$ws2.Range("B3").Value = 389338.4
$ws2.Range("A4").Value = "CFAO MOTORS CI"
$ws2.Range("B4").Value = 365320.25
$ws2.Range("A5").Value = "BRVM - AUTRES SECTEURS"
$ws2.Range("B5").Value = 340394.21
$ws2.Range("A6").Value = "UNIWAX CI"
$ws2.Range("B6").Value = 229856.3
$ws2.Range("A7").Value = "NEI-CEDA CI"
$ws2.Range("B7").Value = 229856.13
$ws2.Range("A8").Value = "SETAO CI"
$ws2.Range("B8").Value = 172656.8
$ws2.Range("A9").Value = "AIR LIQUIDE CI"
$ws2.Range("B9").Value = 168982.88
$ws2.Range("A10").Value = "SUCRIVOIRE"
$ws2.Range("B10").Value = 126443.55
$ws2.Range("A11").Value = "BRVM - DISTRIBUTION"
$ws2.Range("B11").Value = 47856.35
